$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row within the used range (A:R), shifting rows 7-21 down to 8-22,
# without touching the full 16384-column row (keeps dimension at column R).
$ws.Range("A7:R7").Insert(-4121)  # xlShiftDown

# Make sure the date column keeps the same number format as the rest of column D
$ws.Range("D6").Copy()
$ws.Range("D7").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row 7 values
$ws.Cells.Item(7, 1).Value = 10
$ws.Cells.Item(7, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(7, 3).Value = "La Araucanía"
$ws.Cells.Item(7, 4).Value = 44740
$ws.Cells.Item(7, 5).Value = 9
$ws.Cells.Item(7, 6).Value = 100112042
$ws.Cells.Item(7, 7).Value = "Locoto"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 50
$ws.Cells.Item(7, 11).Value = 2500
$ws.Cells.Item(7, 12).Value = 2500
$ws.Cells.Item(7, 13).Value = 2500
$ws.Cells.Item(7, 14).Value = "$/kilo"
$ws.Cells.Item(7, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7, 16).Value = 2500
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Hortaliza"
